# Rename the "filepath" column to "filepaths" and make it hold a
# semicolon-separated list of attachment paths for the first data row
# (David), combining what used to be just "Attachment A.pdf" with
# "Attachment B.pdf".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data cell first, then the header, so that new shared
# strings are appended to the shared-string table in the same order
# a human editing the sheet top-to-bottom-then-header would produce.
$ws.Range("D2").Value = "C:\Users\adavi\OneDrive\Documents\Attachment A.pdf;C:\Users\adavi\OneDrive\Documents\Attachment B.pdf"
$ws.Range("D1").Value = "filepaths"

# Reflect where the user ended up looking after the edit: column D,
# row 1 (the renamed header cell).
$ws.Range("D1").Select()
